$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, pushing existing rows 101:181 down to 102:182
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new weekly record
$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 44574
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100102
$ws.Range("H101").Value = "Cítricos"
$ws.Range("I101").Value = 100102006
$ws.Range("J101").Value = "Pomelo"
$ws.Range("K101").Value = "Start Ruby"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 80
$ws.Range("N101").Value = 14000
$ws.Range("O101").Value = 14000
$ws.Range("P101").Value = 14000
$ws.Range("Q101").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R101").Value = "Región de O'Higgins"
$ws.Range("S101").Value = 1000
$ws.Range("T101").Value = 14
